$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.974.33"
$ws.Range("E2").Value = "  +3.37%  "
$ws.Range("D3").Value = "1.726.35"
$ws.Range("E3").Value = "  +3.03%  "
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.61%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.524"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.38%  "
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "24.11"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +13.57%  "
$ws.Range("E9").Value = "  +3.68%  "
$ws.Range("E10").Value = "  +2.09%  "
$ws.Range("E11").Value = "  +2.17%  "
$ws.Range("D12").Value = "1.970.07"
$ws.Range("E12").Value = "  +3.04%  "
$ws.Range("D13").Value = "1.718.40"
$ws.Range("E13").Value = "  +2.51%  "
$ws.Range("E14").Value = "  +3.68%  "
$ws.Range("E15").Value = "  +5.99%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.88"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.85%  "
$ws.Range("D17").Value = "27.922.28"
$ws.Range("E17").Value = "  +3.22%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "244.08"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.09%  "
$ws.Range("D19").Value = "0.0₃0756"
$ws.Range("E19").Value = "  +2.36%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.90"
$ws.Range("D20").Style = "Normal"
$ws.Range("E21").Value = "  -0.24%  "
$ws.Range("E22").Value = "  +4.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.76"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.86%  "
$ws.Range("E24").Value = "  +0.96%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "149.39"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.00%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.54"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.42%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.83"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.115"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.83%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.19%  "
$ws.Range("E30").Value = "  +2.91%  "
$ws.Range("E31").Value = "  +1.90%  "
$ws.Range("E32").Value = "  +2.87%  "
$ws.Range("E33").Value = "  +3.75%  "
$ws.Range("D34").Value = "1.490.10"
$ws.Range("E34").Value = "  -2.86%  "
$ws.Range("E35").Value = "  -1.56%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.961"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.77%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.612"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.15%  "
$ws.Range("E38").Value = "  +0.75%  "
$ws.Range("E39").Value = "  +0.61%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.07"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.18%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "71.49"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.83%  "
$ws.Range("E42").Value = "  +5.83%  "
$ws.Range("E43").Value = "  -0.20%  "
$ws.Range("B44").Value = "MXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.29"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.27%  "
$ws.Range("B45").Value = "RocketPoolETH"
$ws.Range("C45").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D45").Value = "1.874.15"
$ws.Range("E45").Value = "  +3.09%  "
$ws.Range("E46").Value = "  +1.42%  "
$ws.Range("E47").Value = "  +12.21%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "91.30"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.61%  "
$ws.Range("E49").Value = "  +3.31%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.106"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.60%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.26"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.04%  "
